$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shape = $s.Shapes.Item("Subtitle 2")
$tr = $shape.TextFrame.TextRange

# Find the paragraph that holds the "Section: ..." line (originally
# "Section: L3", a single run) and update the section number to L2.
for ($i = 1; $i -le $tr.Paragraphs().Count; $i++) {
    $para = $tr.Paragraphs($i)
    $text = $para.Text
    if ($text.StartsWith("Section:")) {
        $colonIdx = $text.IndexOf(":")
        # Keep "Section" (everything up to the colon) as-is, and replace
        # the ": L3" suffix with ": L2" - this splits the original single
        # run into "Section" + ": L2".
        $suffix = $para.Characters($colonIdx + 1, $text.Length - $colonIdx)
        $suffix.Text = ": L2"
        break
    }
}
